# Update countries & provincias Spain
#
# The source feed re-ran between two scrapes on 2020-09-29 (04:10 -> 05:27).
# Between those two runs:
#   - A handful of countries' case counts moved just enough to swap places
#     in the (descending, by total cases) ranking with their neighboring row.
#   - Several other countries picked up new case/death numbers without
#     changing rank.
#   - The "last updated" banner text needs bumping to the new time.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Datos actualizados a ..." banner (row 1) -----------------------------
$ws.Range("A1").Value = "Datos actualizados a 29 de Septiembre de 2020 a las 05:27"

# --- Countries that swapped rank with their neighbor -----------------------
# Writing the final text straight into both cells could make the engine
# dedupe the two rows onto a single shared-string slot (since both target
# strings already exist elsewhere in the table at the moment of the write).
# Routing through a unique placeholder first keeps the two swaps as two
# independent writes.

# Rows 22/23: Italia <-> Pakistan
$ws.Range("A22").Value = "__SWAP_PLACEHOLDER_1__"
$ws.Range("A23").Value = "Italia"
$ws.Range("A22").Value = "Pakistan"

# Rows 35/36: Paises Bajos <-> Belgica
$ws.Range("A35").Value = "__SWAP_PLACEHOLDER_2__"
$ws.Range("A36").Value = "Paises Bajos"
$ws.Range("A35").Value = "Belgica"

# Rows 207/208: Timor Oriental <-> Santa Lucia
$ws.Range("A207").Value = "__SWAP_PLACEHOLDER_3__"
$ws.Range("A208").Value = "Timor Oriental"
$ws.Range("A207").Value = "Santa Lucia"

# --- Updated case statistics (columns B:H = Casos totales, Nuevos casos,
#     Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes) ----

# Row 22 - Pakistan
$ws.Range("B22").Value = 311516
$ws.Range("C22").Value = 675
$ws.Range("D22").Value = 296340
$ws.Range("E22").Value = 8702
$ws.Range("G22").Value = 8
$ws.Range("H22").Value = 6474

# Row 23 - Italia
$ws.Range("B23").Value = 311364
$ws.Range("D23").Value = 225190
$ws.Range("E23").Value = 50323
$ws.Range("H23").Value = 35851

# Row 35 - Belgica
$ws.Range("B35").Value = 115353
$ws.Range("C35").Value = 1174
$ws.Range("D35").Value = 19301
$ws.Range("E35").Value = 86065
$ws.Range("G35").Value = 7
$ws.Range("H35").Value = 9987

# Row 36 - Paises Bajos
$ws.Range("B36").Value = 114540
$ws.Range("D36").Value = 0
$ws.Range("E36").Value = 0
$ws.Range("H36").Value = 6380

# Row 39 - Rumania
$ws.Range("B39").Value = 107833
$ws.Range("C39").Value = 58
$ws.Range("D39").Value = 102805
$ws.Range("E39").Value = 3329

# Row 50 - Suecia
$ws.Range("B50").Value = 75537
$ws.Range("C50").Value = 428
$ws.Range("D50").Value = 26957
$ws.Range("E50").Value = 46279
$ws.Range("G50").Value = 12
$ws.Range("H50").Value = 2301

# Row 55 - Nigeria
$ws.Range("B55").Value = 73528
$ws.Range("D55").Value = 63346
$ws.Range("E55").Value = 9568
$ws.Range("H55").Value = 614

# Row 157
$ws.Range("B157").Value = 1891
$ws.Range("C157").Value = 37
$ws.Range("D157").Value = 1208
$ws.Range("E157").Value = 659

# Row 172
$ws.Range("D172").Value = 631
$ws.Range("E172").Value = 46

# Row 173
$ws.Range("B173").Value = 645
$ws.Range("C173").Value = 1
$ws.Range("D173").Value = 549
$ws.Range("E173").Value = 74

# Row 187
$ws.Range("B187").Value = 277
$ws.Range("C187").Value = 1
$ws.Range("E187").Value = 2
